# Loan RBI, Variable Instalments
#
# On the "Repayment schedule" sheet, insert a new (blank) column between the
# existing "In Advance" (M) and "Late" (N) columns. This shifts the old
# N/O/P columns ("Late" / "heading" / "Outstanding") one column to the right
# (to O/P/Q) while leaving their contents/styles intact, and the new column
# N is left blank (only inheriting the header style). Then make
# "Repayment schedule" the active sheet/tab with cell R7 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column at N - shifts old N:P -> O:Q
$ws.Columns("N:N").Insert() | Out-Null

# Match the target column width for the newly inserted column N
$ws.Columns("N:N").ColumnWidth = 10.2

# Make "Repayment schedule" the active sheet/tab and select R7
$ws.Activate() | Out-Null
$ws.Range("R7").Select() | Out-Null
